$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Footnote textbox: merge the three runs that made up footnote "1"
#        into a single run (drop the mid-run formatting split). ---
$footnoteBox = $s.Shapes.Item("TextBox 280")
$tr = $footnoteBox.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Runs in this paragraph (before edit):
#   1: "1"
#   2: " assumes minimal "
#   3: "non-respiratory C "
#   4: "fluxes, …. (Chapin et al. 2006)"
# Work from the tail backwards so earlier-run indices stay valid while we
# fold runs 3 & 4's text into run 2, then blank them out (empty runs are
# dropped on save).
$run4 = $para1.Runs(4, 1)
$run4.Text = ""
$run3 = $para1.Runs(3, 1)
$run3.Text = ""
$run2 = $para1.Runs(2, 1)
$run2.Text = " assumes minimal non-respiratory C fluxes, …. (Chapin et al. 2006)"

# --- 2. Rename the "R_soil" flux label to "R_het_soil" inside the
#        grouped diagram shape. ---
$group = $s.Shapes.Item("Group 323")
$label = $group.GroupItems.Item("Rectangle 322")
$label.TextFrame.TextRange.Text = "▼ R_het_soil"
